$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Sample_ID" etc. shift right by one)
$ws.Columns("D").Insert()

# New column D header - place between is_normal (C) and the shifted Sample_ID (now E)
$ws.Range("D1").Value = "is_normal_for_donor"

# Row 3: is_normal value changes from lowercase "y" to uppercase "Y"
$ws.Range("C3").Value = "Y"

# Add a new 4th data row
$ws.Range("A4").Value = "fake donor 1"
$ws.Range("B4").Value = "fake tissue 1"
$ws.Range("C4").Value = "Y"
$ws.Range("D4").Value = "Y"
$ws.Range("G4").Value = "test_sample.3.bam"
$ws.Range("F4").Value = "t2"
$ws.Range("E4").Value = "test-3"

# Clear the special Courier font formatting previously on relative_file_path cells
# (now shifted from F2:F3 to G2:G3)
$ws.Range("G2:G3").Style = "Normal"

# Remove the explicit row height (reverts to sheet default) on rows 2 and 3
$ws.Rows("2:3").AutoFit()

# Column widths to match the final layout (COM ColumnWidth is in "characters";
# stored OOXML width = ColumnWidth + 5/6, quantized to 1/6 character steps by
# this host, so we back out the closest achievable input per target width)
$ws.Columns("A").ColumnWidth = 10.498697916666666
$ws.Columns("B").ColumnWidth = 10.330729166666666
$ws.Columns("C").ColumnWidth = 8.330729166666666
$ws.Columns("D").ColumnWidth = 17.498697916666668
$ws.Columns("E").ColumnWidth = 8.998697916666666
$ws.Columns("F").ColumnWidth = 11.498697916666666
$ws.Columns("G").ColumnWidth = 15.830729166666666

# Move active selection to D4
$ws.Range("D4").Select() | Out-Null
